# Fruta / hortaliza, semanal
# Insert a new data row before row 89 (pushing the existing rows 89-200 down to 90-201)
# and populate the new row with its own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 89; existing rows shift down automatically.
$ws.Rows(89).Insert()

# Populate the newly inserted row 89 with the new record's data.
$ws.Range("A89").Value = 3
$ws.Range("B89").Value = "Femacal de La Calera"
$ws.Range("C89").Value = "Coquimbo"
$ws.Range("D89").Value = 44483
$ws.Range("E89").Value = 5
$ws.Range("F89").Value = 100112012
$ws.Range("G89").Value = "Espinaca"
$ws.Range("H89").Value = "Sin especificar"
$ws.Range("I89").Value = "Primera"
$ws.Range("J89").Value = 145
$ws.Range("K89").Value = 3000
$ws.Range("L89").Value = 3300
$ws.Range("M89").Value = 3166
$ws.Range("N89").Value = "$/docena de atados (3 kilos)"
$ws.Range("O89").Value = "Provincia de Quillota"
$ws.Range("P89").Value = 1055
$ws.Range("Q89").Value = 3
$ws.Range("R89").Value = "Hortaliza"
